# Apply the edits described by the diff to database_home.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Rename the header of column I from "Internal_links" to "do_not_delete_this_column"
$ws.Range("I1").Value = "do_not_delete_this_column"

# 2. Give I1 the same look as the rest of the header row (A1:H1) instead of
#    the yellow highlight it used to share with J1.
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# 3. Clear the formulas out of I2:I12 (leave the cells blank, but keep their style)
$ws.Range("I2:I12").ClearContents()

# 4. Delete row 13 entirely (it only held the stray I13 cell)
$ws.Rows.Item(13).Delete()

# 5. Match column I's width to column J's width (23.64) now that the long
#    "Internal_links" formula column no longer needs the extra width.
$ws.Range("I1").ColumnWidth = 23.64

# 6. Restore the active selection / view position to match the authored workbook.
$ws.Range("I8").Select()
